# Seguimiento actualizado a 04/05/2015
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Seguimiento")

# --- Fila 14 (registro 9): manuscrito listo para revisión ---
# (se actualiza primero para que la nueva cadena compartida quede antes que
#  "Imágenes corregidas por documentadora" en la tabla de sharedStrings)
$ws.Range("F14").Value = "Manuscrito listo para revisión de María Clemencia"
$ws.Rows.Item(14).RowHeight = 43.5

# --- Fila 6 (registro 1): nueva fecha de entrega y observación actualizada ---
$ws.Range("G6").Value = 42124
$ws.Range("H6").Value = "Imágenes corregidas por documentadora"

# --- Fila 7 (registro 2): nueva fecha de entrega y observación actualizada ---
$ws.Range("G7").Value = 42128
$ws.Range("H7").Value = "Imágenes corregidas por documentadora"

# --- Fila 13 (registro 8): se diligencia entrega y observación ---
# Copiar el formato de fecha (estilo) usado en G6/G7 antes de escribir el valor
$ws.Range("G6").Copy()
$ws.Range("G13").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false
$ws.Range("G13").Value = 42122
$ws.Range("H13").Value = "Imágenes corregidas por documentadora"

# --- Selección activa ---
$ws.Range("G14").Select()
